$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.296.90'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.678.91'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.99'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5340'
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2683'
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06479'
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.93'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07532'
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '1.685.51'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.524'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5780'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008458'
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.80'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '26.335.38'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.904'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.86'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.12'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.208'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.73'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1280'
$ws.Range("E25").Value = '  +6.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.824'
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06498'
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.384'
$ws.Range("E29").Value = '  +3.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.323'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.577'
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.582'
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.032'
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6166'
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.401'
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.244'
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("D39").Value = '1.110.99'
$ws.Range("E39").Value = '  +2.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01619'
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8690'
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.43'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").Value = '1.829.80'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.09'
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("E46").Value = '  -5.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.169'
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4290'
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("E51").Value = '  +1.07%  '
